$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-12-25 Thursday"; new = "2025-12-26 Friday"},
    @{old = "59×62=";  new = "87×13="},
    @{old = "71×59=";  new = "72×51="},
    @{old = "40×18=";  new = "20×90="},
    @{old = "80×36=";  new = "52×13="},
    @{old = "58×21=";  new = "12×28="},
    @{old = "65×64=";  new = "47×44="},
    @{old = "96×64=";  new = "53×78="},
    @{old = "19×21=";  new = "22×19="},
    @{old = "49×79=";  new = "55×35="},
    @{old = "70×78=";  new = "28×91="},
    @{old = "60×80=";  new = "27×12="},
    @{old = "78×14=";  new = "95×26="},
    @{old = "49×27=";  new = "62×16="},
    @{old = "40×57=";  new = "56×62="},
    @{old = "59×89=";  new = "31×79="},
    @{old = "74×59=";  new = "52×29="},
    @{old = "81×26=";  new = "41×56="},
    @{old = "85×84=";  new = "55×29="},
    @{old = "54×65=";  new = "75×22="},
    @{old = "84×76=";  new = "22×25="},
    @{old = "41×59=";  new = "76×11="},
    @{old = "80×93=";  new = "96×57="},
    @{old = "96×44=";  new = "20×25="},
    @{old = "26×54=";  new = "60×24="},
    @{old = "67×70=";  new = "28×98="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
